$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(3, 16).Value = 0
$ws.Cells.Item(3, 17).Value = 0
$ws.Cells.Item(4, 17).Value = 0
$ws.Cells.Item(5, 17).Value = 0
$ws.Cells.Item(6, 17).Value = 0
$ws.Cells.Item(7, 17).Value = 0
$ws.Cells.Item(8, 16).Value = 0
$ws.Cells.Item(8, 17).Value = 0
$ws.Cells.Item(9, 17).Value = 1
$ws.Cells.Item(10, 17).Value = 1
$ws.Cells.Item(11, 17).Value = 1
$ws.Cells.Item(12, 17).Value = 1
$ws.Cells.Item(13, 17).Value = 1
$ws.Cells.Item(14, 17).Value = 1
$ws.Cells.Item(15, 17).Value = 1
$ws.Cells.Item(16, 17).Value = 1
$ws.Cells.Item(17, 17).Value = 2
$ws.Cells.Item(18, 17).Value = 3
$ws.Cells.Item(19, 17).Value = 3
$ws.Cells.Item(20, 17).Value = 4
$ws.Cells.Item(21, 17).Value = 5
$ws.Cells.Item(22, 17).Value = 5
$ws.Cells.Item(23, 17).Value = 6
$ws.Cells.Item(24, 17).Value = 7
$ws.Cells.Item(25, 17).Value = 8
$ws.Cells.Item(26, 17).Value = 8
$ws.Cells.Item(27, 17).Value = 8
$ws.Cells.Item(28, 17).Value = 8
$ws.Cells.Item(29, 17).Value = 8
$ws.Cells.Item(30, 17).Value = 9
$ws.Cells.Item(31, 17).Value = 10
$ws.Cells.Item(32, 17).Value = 11
$ws.Cells.Item(33, 17).Value = 11
$ws.Cells.Item(34, 17).Value = 11
$ws.Cells.Item(35, 17).Value = 11
$ws.Cells.Item(36, 17).Value = 11
$ws.Cells.Item(37, 17).Value = 12
$ws.Cells.Item(38, 17).Value = 12
$ws.Cells.Item(39, 17).Value = 12
$ws.Cells.Item(47, 16).Value = 0
$ws.Cells.Item(47, 17).Value = 3
$ws.Cells.Item(48, 17).Value = 4
$ws.Cells.Item(49, 17).Value = 4
$ws.Cells.Item(50, 17).Value = 4
$ws.Cells.Item(51, 17).Value = 4
$ws.Cells.Item(52, 17).Value = 4
$ws.Cells.Item(53, 15).Value = '[''Ireland'', ''Colombia'', ''Costa Rica'', ''Argentina'']'
$ws.Cells.Item(53, 17).Value = 4
$ws.Cells.Item(54, 15).Value = '[''Ireland'', ''Colombia'', ''Costa Rica'', ''Argentina'']'
$ws.Cells.Item(54, 17).Value = 4
$ws.Cells.Item(55, 15).Value = '[''Ireland'', ''Colombia'', ''Costa Rica'', ''Argentina'']'
$ws.Cells.Item(55, 17).Value = 4
$ws.Cells.Item(56, 15).Value = '[''Ireland'', ''Colombia'', ''Costa Rica'', ''Argentina'']'
$ws.Cells.Item(56, 17).Value = 4
$ws.Cells.Item(57, 15).Value = '[''Ireland'', ''Colombia'', ''Costa Rica'', ''Argentina'']'
$ws.Cells.Item(57, 17).Value = 5
$ws.Cells.Item(58, 15).Value = '[''Ireland'', ''Colombia'', ''Costa Rica'', ''Argentina'']'
$ws.Cells.Item(58, 17).Value = 5
$ws.Cells.Item(59, 15).Value = '[''Ireland'', ''Colombia'', ''Costa Rica'', ''Argentina'']'
$ws.Cells.Item(59, 17).Value = 6
$ws.Cells.Item(60, 15).Value = '[''Colombia'', ''Argentina'']'
$ws.Cells.Item(60, 17).Value = 6
$ws.Cells.Item(61, 15).Value = '[''Colombia'', ''Argentina'']'
$ws.Cells.Item(61, 17).Value = 7
$ws.Cells.Item(62, 15).Value = '[''Colombia'', ''Argentina'']'
$ws.Cells.Item(62, 17).Value = 8
$ws.Cells.Item(63, 15).Value = '[''Colombia'', ''Scotland'', ''Austria'', ''Argentina'']'
$ws.Cells.Item(63, 17).Value = 9
$ws.Cells.Item(64, 15).Value = '[''Colombia'', ''Scotland'', ''Austria'', ''Argentina'']'
$ws.Cells.Item(64, 17).Value = 9
$ws.Cells.Item(65, 15).Value = '[''Colombia'', ''Scotland'', ''Austria'', ''Argentina'']'
$ws.Cells.Item(65, 17).Value = 9
$ws.Cells.Item(66, 15).Value = '[''Colombia'', ''Scotland'', ''Austria'', ''Argentina'']'
$ws.Cells.Item(66, 17).Value = 9
$ws.Cells.Item(67, 15).Value = '[''Colombia'', ''Scotland'', ''Austria'', ''Argentina'']'
$ws.Cells.Item(67, 17).Value = 9
$ws.Cells.Item(68, 15).Value = '[''Colombia'', ''Scotland'', ''Austria'', ''Argentina'']'
$ws.Cells.Item(68, 17).Value = 9
$ws.Cells.Item(69, 15).Value = '[''Colombia'', ''Scotland'', ''Austria'', ''Argentina'']'
$ws.Cells.Item(69, 17).Value = 10
$ws.Cells.Item(70, 15).Value = '[''Colombia'', ''Scotland'', ''Austria'', ''Argentina'']'
$ws.Cells.Item(70, 16).Value = 0
$ws.Cells.Item(70, 17).Value = 10
$ws.Cells.Item(71, 15).Value = '[''Colombia'', ''Scotland'', ''Austria'', ''Argentina'']'
$ws.Cells.Item(71, 17).Value = 11
$ws.Cells.Item(72, 15).Value = '[''Colombia'', ''Scotland'', ''Austria'', ''Argentina'']'
$ws.Cells.Item(72, 17).Value = 12
$ws.Cells.Item(73, 15).Value = '[''Colombia'', ''Scotland'', ''Austria'', ''Argentina'']'
$ws.Cells.Item(73, 17).Value = 13
$ws.Cells.Item(88, 16).Value = 0
$ws.Cells.Item(88, 17).Value = 3
$ws.Cells.Item(89, 17).Value = 4
$ws.Cells.Item(90, 17).Value = 5
$ws.Cells.Item(91, 17).Value = 5
$ws.Cells.Item(92, 17).Value = 5
$ws.Cells.Item(93, 17).Value = 5
$ws.Cells.Item(94, 17).Value = 5
$ws.Cells.Item(95, 17).Value = 5
$ws.Cells.Item(96, 17).Value = 5
$ws.Cells.Item(97, 17).Value = 5
$ws.Cells.Item(98, 17).Value = 5
$ws.Cells.Item(99, 17).Value = 6
$ws.Cells.Item(100, 16).Value = 0
$ws.Cells.Item(100, 17).Value = 6
$ws.Cells.Item(101, 15).Value = '[''Netherlands'', ''Italy'']'
$ws.Cells.Item(101, 17).Value = 7
$ws.Cells.Item(102, 15).Value = '[''Netherlands'', ''Italy'']'
$ws.Cells.Item(102, 17).Value = 7
$ws.Cells.Item(103, 17).Value = 8
$ws.Cells.Item(104, 17).Value = 9
$ws.Cells.Item(105, 17).Value = 10
$ws.Cells.Item(106, 17).Value = 11
$ws.Cells.Item(107, 17).Value = 11
$ws.Cells.Item(108, 16).Value = 0
$ws.Cells.Item(108, 17).Value = 11
$ws.Cells.Item(109, 17).Value = 12
$ws.Cells.Item(110, 17).Value = 12
